# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(0.0003714022599530242, 0.05231270169004087, 3.082599426703578, 6.48142807727062, 9.616711607924191)
    3 = @(0.02258322285507441, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 2.888439239842931)
    4 = @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 2.652525447291612)
    5 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    6 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
